$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.413.40'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.519.85'
$ws.Range("E3").Value = '  +1.23%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.57'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.46'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.543.52'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.965.37'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.56'
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.323.48'
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.532.23'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.29'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.61'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.90'
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.22'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.434'
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.994'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0789'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.83'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.74'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("E32").Value = '  -5.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.49'
$ws.Range("E33").Value = '  +3.35%  '
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.90'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.39'
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("E38").Value = '  -6.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.06'
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("E40").Value = '  -5.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.837'
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '296.72'
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.70'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.603'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.82'
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0938'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.86'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.32'
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("E51").Value = '  -3.71%  '
